$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 (I0) and J1 (IF), matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Row-by-row data for the new I0 / IF columns (rows 2-53)
$data = @{
    2 = @(9, 9)
    3 = @(9, 9)
    4 = @(8, 8)
    5 = @(7, 7)
    6 = @(9, 9)
    7 = @(6, 6)
    8 = @(8, 8)
    9 = @(6, 6)
    10 = @(7, 7)
    11 = @(1, 1)
    12 = @(11, 12)
    13 = @(6, 6)
    14 = @(6, 6)
    15 = @(6, 6)
    16 = @(7, 7)
    17 = @(8, 8)
    18 = @(7, 8)
    19 = @(8, 8)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(8, 8)
    26 = @(6, 6)
    27 = @(6, 6)
    28 = @(5, 6)
    29 = @(6, 6)
    30 = @(8, 8)
    31 = @(8, 8)
    32 = @(6, 6)
    33 = @(7, 7)
    34 = @(7, 7)
    35 = @(9, 9)
    36 = @(8, 9)
    37 = @(9, 9)
    38 = @(7, 8)
    39 = @(7, 7)
    40 = @(7, 7)
    41 = @(7, 7)
    42 = @(9, 9)
    43 = @(6, 6)
    44 = @(9, 9)
    45 = @(6, 7)
    46 = @(8, 9)
    47 = @(6, 7)
    48 = @(7, 7)
    49 = @(5, 5)
    50 = @(7, 7)
    51 = @(5, 5)
    52 = @(8, 8)
    53 = @(2, 2)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Host "Done updating I0/IF columns"
